{"js": "// Replace each arithmetic expression's old text with its new text.\n// Every \"old\" value below occurs exactly once in the document, so a\n// matchCase, non-partial-word search safely finds the single target run.\nconst replacements = [\n  [\n    \"82-59=23\",\n    \"32+19=51\"\n  ],\n  [\n    \"8+73=81\",\n    \"8+91=99\"\n  ],\n  [\n    \"70-42=28\",\n    \"1+51=52\"\n  ],\n  [\n    \"28+63=91\",\n    \"93-9=84\"\n  ],\n  [\n    \"68-16=52\",\n    \"91+3=94\"\n  ],\n  [\n    \"11+6=17\",\n    \"86-44=42\"\n  ],\n  [\n    \"28-3=25\",\n    \"29+18=47\"\n  ],\n  [\n    \"9+15=24\",\n    \"43-0=43\"\n  ],\n  [\n    \"5+51=56\",\n    \"48-21=27\"\n  ],\n  [\n    \"20+1=21\",\n    \"76-19=57\"\n  ],\n  [\n    \"71+26=97\",\n    \"89-62=27\"\n  ],\n  [\n    \"29+38=67\",\n    \"39-26=13\"\n  ],\n  [\n    \"59-45=14\",\n    \"17+54=71\"\n  ],\n  [\n    \"38+47=85\",\n    \"79-67=12\"\n  ],\n  [\n    \"4+57=61\",\n    \"81-51=30\"\n  ],\n  [\n    \"12+38=50\",\n    \"9+39=48\"\n  ],\n  [\n    \"81-4=77\",\n    \"37+56=93\"\n  ],\n  [\n    \"68+9=77\",\n    \"0+56=56\"\n  ],\n  [\n    \"74-36=38\",\n    \"8+21=29\"\n  ],\n  [\n    \"54-25=29\",\n    \"80+5=85\"\n  ],\n  [\n    \"85-33=52\",\n    \"17+59=76\"\n  ],\n  [\n    \"66-56=10\",\n    \"83-66=17\"\n  ],\n  [\n    \"34+43=77\",\n    \"62+15=77\"\n  ],\n  [\n    \"73-63=10\",\n    \"1+74=75\"\n  ],\n  [\n    \"13+13=26\",\n    \"55+8=63\"\n  ],\n  [\n    \"22+29=51\",\n    \"26+20=46\"\n  ],\n  [\n    \"17+76=93\",\n    \"61-29=32\"\n  ],\n  [\n    \"3+2=5\",\n    \"86-85=1\"\n  ],\n  [\n    \"57-48=9\",\n    \"68-35=33\"\n  ],\n  [\n    \"67-6=61\",\n    \"80-46=34\"\n  ],\n  [\n    \"8+64=72\",\n    \"49+31=80\"\n  ],\n  [\n    \"0+11=11\",\n    \"51-40=11\"\n  ],\n  [\n    \"68+19=87\",\n    \"36-1=35\"\n  ],\n  [\n    \"21+26=47\",\n    \"8-6=2\"\n  ],\n  [\n    \"91-11=80\",\n    \"51-6=45\"\n  ],\n  [\n    \"78-45=33\",\n    \"20+76=96\"\n  ],\n  [\n    \"1+24=25\",\n    \"4+7=11\"\n  ],\n  [\n    \"37+37=74\",\n    \"67+7=74\"\n  ],\n  [\n    \"14+44=58\",\n    \"1+91=92\"\n  ],\n  [\n    \"3+74=77\",\n    \"57+29=86\"\n  ],\n  [\n    \"79-36=43\",\n    \"74-70=4\"\n  ],\n  [\n    \"36+37=73\",\n    \"15+60=75\"\n  ],\n  [\n    \"30+26=56\",\n    \"4-2=2\"\n  ],\n  [\n    \"35+61=96\",\n    \"99-6=93\"\n  ],\n  [\n    \"64-48=16\",\n    \"51+21=72\"\n  ],\n  [\n    \"37+7=44\",\n    \"47+17=64\"\n  ],\n  [\n    \"13+11=24\",\n    \"37+45=82\"\n  ],\n  [\n    \"82+7=89\",\n    \"6+68=74\"\n  ],\n  [\n    \"10+55=65\",\n    \"87-66=21\"\n  ],\n  [\n    \"2+18=20\",\n    \"27-24=3\"\n  ],\n  [\n    \"30-20=10\",\n    \"83-65=18\"\n  ],\n  [\n    \"0+6=6\",\n    \"45+53=98\"\n  ],\n  [\n    \"89-44=45\",\n    \"71+18=89\"\n  ],\n  [\n    \"60-56=4\",\n    \"63-12=51\"\n  ],\n  [\n    \"96-43=53\",\n    \"25+63=88\"\n  ],\n  [\n    \"35+3=38\",\n    \"47+21=68\"\n  ],\n  [\n    \"6+4=10\",\n    \"37+30=67\"\n  ],\n  [\n    \"86-57=29\",\n    \"15-2=13\"\n  ],\n  [\n    \"57-10=47\",\n    \"80-33=47\"\n  ],\n  [\n    \"13+2=15\",\n    \"92-35=57\"\n  ],\n  [\n    \"68-9=59\",\n    \"6+56=62\"\n  ],\n  [\n    \"85-77=8\",\n    \"20+36=56\"\n  ],\n  [\n    \"36+59=95\",\n    \"67+5=72\"\n  ],\n  [\n    \"2+39=41\",\n    \"32-1=31\"\n  ],\n  [\n    \"94-2=92\",\n    \"33+38=71\"\n  ],\n  [\n    \"61-11=50\",\n    \"13+47=60\"\n  ],\n  [\n    \"53+35=88\",\n    \"22+52=74\"\n  ],\n  [\n    \"89-37=52\",\n    \"4+81=85\"\n  ],\n  [\n    \"34+0=34\",\n    \"29-25=4\"\n  ],\n  [\n    \"30+12=42\",\n    \"46+21=67\"\n  ],\n  [\n    \"18+18=36\",\n    \"92+7=99\"\n  ],\n  [\n    \"89-79=10\",\n    \"57-28=29\"\n  ],\n  [\n    \"35+49=84\",\n    \"88-51=37\"\n  ],\n  [\n    \"82+13=95\",\n    \"54-41=13\"\n  ],\n  [\n    \"49+39=88\",\n    \"81-35=46\"\n  ],\n  [\n    \"27+4=31\",\n    \"20+52=72\"\n  ],\n  [\n    \"65-4=61\",\n    \"0+78=78\"\n  ],\n  [\n    \"71+27=98\",\n    \"88-26=62\"\n  ],\n  [\n    \"34+62=96\",\n    \"8-4=4\"\n  ],\n  [\n    \"97-27=70\",\n    \"48+2=50\"\n  ],\n  [\n    \"47-15=32\",\n    \"8+16=24\"\n  ],\n  [\n    \"34+53=87\",\n    \"49-7=42\"\n  ],\n  [\n    \"34-7=27\",\n    \"13+70=83\"\n  ],\n  [\n    \"61-5=56\",\n    \"99-91=8\"\n  ],\n  [\n    \"43-27=16\",\n    \"78-74=4\"\n  ],\n  [\n    \"75-11=64\",\n    \"54+42=96\"\n  ],\n  [\n    \"20+30=50\",\n    \"52+0=52\"\n  ],\n  [\n    \"43+16=59\",\n    \"5+40=45\"\n  ],\n  [\n    \"98-15=83\",\n    \"85-65=20\"\n  ],\n  [\n    \"97-85=12\",\n    \"37-36=1\"\n  ],\n  [\n    \"68-61=7\",\n    \"18-2=16\"\n  ],\n  [\n    \"74-5=69\",\n    \"27+34=61\"\n  ],\n  [\n    \"32+24=56\",\n    \"27-21=6\"\n  ],\n  [\n    \"89+0=89\",\n    \"65-33=32\"\n  ],\n  [\n    \"29+40=69\",\n    \"59+9=68\"\n  ],\n  [\n    \"24-22=2\",\n    \"23+23=46\"\n  ],\n  [\n    \"43+3=46\",\n    \"47-28=19\"\n  ],\n  [\n    \"95-12=83\",\n    \"86-36=50\"\n  ],\n  [\n    \"71-56=15\",\n    \"5+62=67\"\n  ],\n  [\n    \"59-6=53\",\n    \"26-12=14\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic expression's old text with its new text.\n# Every \"old\" value below occurs exactly once in the document, so\n# Find/Replace (non-wildcard, whole document) safely targets only\n# the single matching run each time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"82-59=23\", \"32+19=51\"),\n  @(\"8+73=81\", \"8+91=99\"),\n  @(\"70-42=28\", \"1+51=52\"),\n  @(\"28+63=91\", \"93-9=84\"),\n  @(\"68-16=52\", \"91+3=94\"),\n  @(\"11+6=17\", \"86-44=42\"),\n  @(\"28-3=25\", \"29+18=47\"),\n  @(\"9+15=24\", \"43-0=43\"),\n  @(\"5+51=56\", \"48-21=27\"),\n  @(\"20+1=21\", \"76-19=57\"),\n  @(\"71+26=97\", \"89-62=27\"),\n  @(\"29+38=67\", \"39-26=13\"),\n  @(\"59-45=14\", \"17+54=71\"),\n  @(\"38+47=85\", \"79-67=12\"),\n  @(\"4+57=61\", \"81-51=30\"),\n  @(\"12+38=50\", \"9+39=48\"),\n  @(\"81-4=77\", \"37+56=93\"),\n  @(\"68+9=77\", \"0+56=56\"),\n  @(\"74-36=38\", \"8+21=29\"),\n  @(\"54-25=29\", \"80+5=85\"),\n  @(\"85-33=52\", \"17+59=76\"),\n  @(\"66-56=10\", \"83-66=17\"),\n  @(\"34+43=77\", \"62+15=77\"),\n  @(\"73-63=10\", \"1+74=75\"),\n  @(\"13+13=26\", \"55+8=63\"),\n  @(\"22+29=51\", \"26+20=46\"),\n  @(\"17+76=93\", \"61-29=32\"),\n  @(\"3+2=5\", \"86-85=1\"),\n  @(\"57-48=9\", \"68-35=33\"),\n  @(\"67-6=61\", \"80-46=34\"),\n  @(\"8+64=72\", \"49+31=80\"),\n  @(\"0+11=11\", \"51-40=11\"),\n  @(\"68+19=87\", \"36-1=35\"),\n  @(\"21+26=47\", \"8-6=2\"),\n  @(\"91-11=80\", \"51-6=45\"),\n  @(\"78-45=33\", \"20+76=96\"),\n  @(\"1+24=25\", \"4+7=11\"),\n  @(\"37+37=74\", \"67+7=74\"),\n  @(\"14+44=58\", \"1+91=92\"),\n  @(\"3+74=77\", \"57+29=86\"),\n  @(\"79-36=43\", \"74-70=4\"),\n  @(\"36+37=73\", \"15+60=75\"),\n  @(\"30+26=56\", \"4-2=2\"),\n  @(\"35+61=96\", \"99-6=93\"),\n  @(\"64-48=16\", \"51+21=72\"),\n  @(\"37+7=44\", \"47+17=64\"),\n  @(\"13+11=24\", \"37+45=82\"),\n  @(\"82+7=89\", \"6+68=74\"),\n  @(\"10+55=65\", \"87-66=21\"),\n  @(\"2+18=20\", \"27-24=3\"),\n  @(\"30-20=10\", \"83-65=18\"),\n  @(\"0+6=6\", \"45+53=98\"),\n  @(\"89-44=45\", \"71+18=89\"),\n  @(\"60-56=4\", \"63-12=51\"),\n  @(\"96-43=53\", \"25+63=88\"),\n  @(\"35+3=38\", \"47+21=68\"),\n  @(\"6+4=10\", \"37+30=67\"),\n  @(\"86-57=29\", \"15-2=13\"),\n  @(\"57-10=47\", \"80-33=47\"),\n  @(\"13+2=15\", \"92-35=57\"),\n  @(\"68-9=59\", \"6+56=62\"),\n  @(\"85-77=8\", \"20+36=56\"),\n  @(\"36+59=95\", \"67+5=72\"),\n  @(\"2+39=41\", \"32-1=31\"),\n  @(\"94-2=92\", \"33+38=71\"),\n  @(\"61-11=50\", \"13+47=60\"),\n  @(\"53+35=88\", \"22+52=74\"),\n  @(\"89-37=52\", \"4+81=85\"),\n  @(\"34+0=34\", \"29-25=4\"),\n  @(\"30+12=42\", \"46+21=67\"),\n  @(\"18+18=36\", \"92+7=99\"),\n  @(\"89-79=10\", \"57-28=29\"),\n  @(\"35+49=84\", \"88-51=37\"),\n  @(\"82+13=95\", \"54-41=13\"),\n  @(\"49+39=88\", \"81-35=46\"),\n  @(\"27+4=31\", \"20+52=72\"),\n  @(\"65-4=61\", \"0+78=78\"),\n  @(\"71+27=98\", \"88-26=62\"),\n  @(\"34+62=96\", \"8-4=4\"),\n  @(\"97-27=70\", \"48+2=50\"),\n  @(\"47-15=32\", \"8+16=24\"),\n  @(\"34+53=87\", \"49-7=42\"),\n  @(\"34-7=27\", \"13+70=83\"),\n  @(\"61-5=56\", \"99-91=8\"),\n  @(\"43-27=16\", \"78-74=4\"),\n  @(\"75-11=64\", \"54+42=96\"),\n  @(\"20+30=50\", \"52+0=52\"),\n  @(\"43+16=59\", \"5+40=45\"),\n  @(\"98-15=83\", \"85-65=20\"),\n  @(\"97-85=12\", \"37-36=1\"),\n  @(\"68-61=7\", \"18-2=16\"),\n  @(\"74-5=69\", \"27+34=61\"),\n  @(\"32+24=56\", \"27-21=6\"),\n  @(\"89+0=89\", \"65-33=32\"),\n  @(\"29+40=69\", \"59+9=68\"),\n  @(\"24-22=2\", \"23+23=46\"),\n  @(\"43+3=46\", \"47-28=19\"),\n  @(\"95-12=83\", \"86-36=50\"),\n  @(\"71-56=15\", \"5+62=67\"),\n  @(\"59-6=53\", \"26-12=14\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n"}
